# daily auto push: 2025-10-03 09:26 UTC
# Append a new data row (row 56) to Sheet1 with the day's latest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 56

# Column A holds a date-like string ("2025/10/03"); force text formatting
# first so Excel doesn't auto-convert it into a date serial number, then
# clear the format again so the cell keeps the workbook's default style
# (matching the other data rows, which carry no explicit style index).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/03"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "金"
$ws.Cells.Item($row, 3).Value = 17
$ws.Cells.Item($row, 4).Value = 201
